$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 236, shifting existing rows 236:281 down to 238:283.
$ws.Rows.Item(236).Resize(2).EntireRow.Insert()

# New row 236 - "Primera" quality for the new weekly report date (2021-12-13 -> 44543)
$ws.Cells.Item(236, 1).Value = 8
$ws.Cells.Item(236, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44543
$ws.Cells.Item(236, 5).Value = 4
$ws.Cells.Item(236, 6).Value = 100112009
$ws.Cells.Item(236, 7).Value = "Acelga"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 2200
$ws.Cells.Item(236, 11).Value = 550
$ws.Cells.Item(236, 12).Value = 600
$ws.Cells.Item(236, 13).Value = 575
$ws.Cells.Item(236, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(236, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(236, 16).Value = 288
$ws.Cells.Item(236, 17).Value = 2
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# New row 237 - "Segunda" quality for the same new weekly report date
$ws.Cells.Item(237, 1).Value = 8
$ws.Cells.Item(237, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44543
$ws.Cells.Item(237, 5).Value = 4
$ws.Cells.Item(237, 6).Value = 100112009
$ws.Cells.Item(237, 7).Value = "Acelga"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Segunda"
$ws.Cells.Item(237, 10).Value = 1500
$ws.Cells.Item(237, 11).Value = 450
$ws.Cells.Item(237, 12).Value = 500
$ws.Cells.Item(237, 13).Value = 475
$ws.Cells.Item(237, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(237, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(237, 16).Value = 238
$ws.Cells.Item(237, 17).Value = 2
$ws.Cells.Item(237, 18).Value = "Hortaliza"
